# Sun, May 03, 2020  9:08:22 PM
#
# The deck's Slide Master currently carries the "Integral" design theme
# (ppt/theme/theme1.xml) while the Notes Master carries the default
# "Office Theme" colours (ppt/theme/theme2.xml). The author switched the
# Slide Master's theme colours over to the stock "Office Theme" palette
# (the built-in blue/orange Office colour scheme), the same twelve
# theme colours that the Notes Master's theme already used.
#
# Re-colour every slot of the Slide Master's theme colour scheme
# (Background/Text 1-2, Accent 1-6, Hyperlink, Followed Hyperlink) to
# the standard Office Theme RGB values.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

function Set-ThemeRGB($index, $r, $g, $b) {
    $themeColors.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# msoThemeDark1 / Background 1
Set-ThemeRGB 1  0x00 0x00 0x00
# msoThemeLight1 / Text 1
Set-ThemeRGB 2  0xFF 0xFF 0xFF
# msoThemeDark2 / Background 2
Set-ThemeRGB 3  0x44 0x54 0x6A
# msoThemeLight2 / Text 2
Set-ThemeRGB 4  0xE7 0xE6 0xE6
# Accent 1
Set-ThemeRGB 5  0x5B 0x9B 0xD5
# Accent 2
Set-ThemeRGB 6  0xED 0x7D 0x31
# Accent 3
Set-ThemeRGB 7  0xA5 0xA5 0xA5
# Accent 4
Set-ThemeRGB 8  0xFF 0xC0 0x00
# Accent 5
Set-ThemeRGB 9  0x44 0x72 0xC4
# Accent 6
Set-ThemeRGB 10 0x70 0xAD 0x47
# Hyperlink
Set-ThemeRGB 11 0x05 0x63 0xC1
# Followed Hyperlink
Set-ThemeRGB 12 0x95 0x4F 0x72
